$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "24.385.10"
Set-TextValue $ws.Range("E2") "  -1.56%  "

Set-TextValue $ws.Range("D3") "1.653.61"
Set-TextValue $ws.Range("E3") "  -2.77%  "

Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  +0.08%  "

Set-TextValue $ws.Range("D5") "311.34"
Set-TextValue $ws.Range("E5") "  -0.98%  "

Set-TextValue $ws.Range("D7") "0.3914"
Set-TextValue $ws.Range("E7") "  -1.61%  "

Set-TextValue $ws.Range("D8") "0.3903"
Set-TextValue $ws.Range("E8") "  -3.42%  "

Set-TextValue $ws.Range("D9") "1.001"
Set-TextValue $ws.Range("E9") "  -0.26%  "

Set-TextValue $ws.Range("D10") "1.380"
Set-TextValue $ws.Range("E10") "  -5.72%  "

Set-TextValue $ws.Range("D11") "50.01"
Set-TextValue $ws.Range("E11") "  -6.60%  "

Set-TextValue $ws.Range("D12") "0.08548"
Set-TextValue $ws.Range("E12") "  -2.80%  "

Set-TextValue $ws.Range("D13") "24.95"
Set-TextValue $ws.Range("E13") "  -5.07%  "

Set-TextValue $ws.Range("D14") "7.208"
Set-TextValue $ws.Range("E14") "  -4.05%  "

Set-TextValue $ws.Range("D15") "0.00001304"
Set-TextValue $ws.Range("E15") "  -2.74%  "

Set-TextValue $ws.Range("D16") "7.605"
Set-TextValue $ws.Range("E16") "  -4.59%  "

Set-TextValue $ws.Range("D17") "1.652.24"
Set-TextValue $ws.Range("E17") "  -2.42%  "

Set-TextValue $ws.Range("D18") "93.17"
Set-TextValue $ws.Range("E18") "  -2.45%  "

Set-TextValue $ws.Range("D19") "0.06948"
Set-TextValue $ws.Range("E19") "  -3.19%  "

Set-TextValue $ws.Range("D20") "21.01"
Set-TextValue $ws.Range("E20") "  +0.94%  "

Set-TextValue $ws.Range("D21") "7.004"
Set-TextValue $ws.Range("E21") "  -4.31%  "

Set-TextValue $ws.Range("D22") "1.003"
Set-TextValue $ws.Range("E22") "  +0.06%  "

Set-TextValue $ws.Range("E23") "  -4.01%  "

Set-TextValue $ws.Range("D24") "24.402.27"
Set-TextValue $ws.Range("E24") "  -1.47%  "

Set-TextValue $ws.Range("D25") "2.340"
Set-TextValue $ws.Range("E25") "  -1.77%  "

Set-TextValue $ws.Range("D26") "2.788"
Set-TextValue $ws.Range("E26") "  -3.66%  "

Set-TextValue $ws.Range("D27") "22.65"
Set-TextValue $ws.Range("E27") "  -1.83%  "

Set-TextValue $ws.Range("D28") "158.73"
Set-TextValue $ws.Range("E28") "  -2.10%  "

Set-TextValue $ws.Range("D29") "5.721"
Set-TextValue $ws.Range("E29") "  -6.60%  "

Set-TextValue $ws.Range("D30") "145.14"
Set-TextValue $ws.Range("E30") "  +0.92%  "

Set-TextValue $ws.Range("D31") "8.222"
Set-TextValue $ws.Range("E31") "  -0.57%  "

Set-TextValue $ws.Range("D32") "2.507"
Set-TextValue $ws.Range("E32") "  +10.86%  "

Set-TextValue $ws.Range("D33") "1.838.19"
Set-TextValue $ws.Range("E33") "  -6.23%  "

Set-TextValue $ws.Range("D34") "0.03015"
Set-TextValue $ws.Range("E34") "  -4.62%  "

Set-TextValue $ws.Range("D35") "0.08110"
Set-TextValue $ws.Range("E35") "  -5.57%  "

Set-TextValue $ws.Range("D36") "0.9973"
Set-TextValue $ws.Range("E36") "  -2.65%  "

Set-TextValue $ws.Range("D37") "6.854"
Set-TextValue $ws.Range("E37") "  -6.29%  "

Set-TextValue $ws.Range("D38") "0.2760"
Set-TextValue $ws.Range("E38") "  -2.94%  "

Set-TextValue $ws.Range("D39") "0.09471"
Set-TextValue $ws.Range("E39") "  +0.48%  "

Set-TextValue $ws.Range("D40") "1.493"
Set-TextValue $ws.Range("E40") "  +1.14%  "

Set-TextValue $ws.Range("D41") "10.19"
Set-TextValue $ws.Range("E41") "  -4.64%  "

Set-TextValue $ws.Range("D42") "0.7791"
Set-TextValue $ws.Range("E42") "  -5.94%  "

Set-TextValue $ws.Range("D43") "13.33"
Set-TextValue $ws.Range("E43") "  -5.78%  "

Set-TextValue $ws.Range("D44") "16.30"
Set-TextValue $ws.Range("E44") "  -7.20%  "

Set-TextValue $ws.Range("D45") "2.551"
Set-TextValue $ws.Range("E45") "  -5.42%  "

Set-TextValue $ws.Range("D46") "0.7011"
Set-TextValue $ws.Range("E46") "  -5.34%  "

Set-TextValue $ws.Range("E47") "  -1.89%  "

Set-TextValue $ws.Range("B48") "Cronos"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.08587"
Set-TextValue $ws.Range("E48") "  +2.49%  "

Set-TextValue $ws.Range("B49") "Frax"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D49") "1.002"
Set-TextValue $ws.Range("E49") "  -0.08%  "

Set-TextValue $ws.Range("D50") "1.302"
Set-TextValue $ws.Range("E50") "  -5.21%  "

Set-TextValue $ws.Range("D51") "136.32"
Set-TextValue $ws.Range("E51") "  -2.04%  "
